# Refresh the crypto price/volume table (cryptos.xlsx) with the latest
# scrape: most rows keep their rank but get new Price (D) / Volume(1h) (E)
# figures, while a couple of coins swapped rank with their neighbour
# (Litecoin <-> InternetComputer(DFINITY) at rows 22/23, and
# Stellar <-> Kaspa at rows 38/39), which also moves their Coin name (B)
# and Link (C) text onto the other row.
#
# All cells in this sheet are stored as literal text (even things like
# "1.00" or "86.98"), so any column-D price that reads as a plain decimal
# number is written with a leading apostrophe. That is the standard
# Excel-UI way to force text entry, and it keeps those values from being
# silently reinterpreted as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index (B=2, C=3, D=4, E=5), new value.
$updates = @(
    @(2, 4, "40.131.62"),
    @(2, 5, "  +2.64%  "),
    @(3, 4, "2.244.30"),
    @(3, 5, "  +0.41%  "),
    @(4, 5, "  +0.02%  "),
    @(5, 4, "'294.91"),
    @(5, 5, "  -1.11%  "),
    @(6, 4, "'86.98"),
    @(6, 5, "  +6.64%  "),
    @(7, 5, "  +0.80%  "),
    @(8, 5, "  +0.06%  "),
    @(9, 4, "'0.474"),
    @(9, 5, "  +2.78%  "),
    @(10, 4, "'31.11"),
    @(10, 5, "  +11.15%  "),
    @(11, 4, "'0.0802"),
    @(11, 5, "  +3.37%  "),
    @(12, 4, "'47.11"),
    @(12, 5, "  +2.48%  "),
    @(13, 5, "  +0.51%  "),
    @(14, 4, "'6.46"),
    @(14, 5, "  +5.68%  "),
    @(15, 4, "2.587.79"),
    @(15, 5, "  +0.79%  "),
    @(16, 4, "'14.27"),
    @(16, 5, "  +1.28%  "),
    @(17, 4, "2.252.66"),
    @(17, 5, "  +0.98%  "),
    @(18, 4, "'0.734"),
    @(18, 5, "  +2.31%  "),
    @(19, 4, "40.034.45"),
    @(19, 5, "  +2.65%  "),
    @(20, 4, "0.0₃0895"),
    @(20, 5, "  +3.68%  "),
    @(21, 4, "'5.88"),
    @(21, 5, "  +2.24%  "),
    @(22, 2, "InternetComputer(DFINITY)"),
    @(22, 3, "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"),
    @(22, 4, "'10.64"),
    @(22, 5, "  +7.12%  "),
    @(23, 2, "Litecoin"),
    @(23, 3, "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"),
    @(23, 4, "'65.82"),
    @(23, 5, "  +0.94%  "),
    @(24, 4, "'236.69"),
    @(24, 5, "  +4.78%  "),
    @(25, 5, "  +0.01%  "),
    @(26, 4, "'2.47"),
    @(26, 5, "  +3.40%  "),
    @(27, 4, "'1.85"),
    @(27, 5, "  +6.57%  "),
    @(28, 4, "'23.13"),
    @(28, 5, "  +3.62%  "),
    @(29, 4, "'2.23"),
    @(29, 5, "  +1.89%  "),
    @(30, 4, "'9.31"),
    @(30, 5, "  +4.34%  "),
    @(31, 4, "'34.22"),
    @(31, 5, "  +7.62%  "),
    @(32, 4, "'152.77"),
    @(32, 5, "  +2.55%  "),
    @(33, 5, "  -0.02%  "),
    @(34, 4, "'4.90"),
    @(34, 5, "  +2.19%  "),
    @(35, 4, "'0.0718"),
    @(35, 5, "  +4.63%  "),
    @(36, 4, "'2.39"),
    @(36, 5, "  +2.41%  "),
    @(37, 4, "'16.72"),
    @(37, 5, "  +13.00%  "),
    @(38, 2, "Kaspa"),
    @(38, 3, "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @(38, 4, "'0.102"),
    @(38, 5, "  +5.46%  "),
    @(39, 2, "Stellar"),
    @(39, 3, "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @(39, 4, "'0.112"),
    @(39, 5, "  +2.30%  "),
    @(40, 4, "'2.72"),
    @(40, 5, "  +2.42%  "),
    @(41, 4, "'1.69"),
    @(41, 5, "  +4.80%  "),
    @(42, 4, "'3.82"),
    @(42, 5, "  +4.45%  "),
    @(43, 4, "1.990.48"),
    @(43, 5, "  +4.24%  "),
    @(44, 4, "'2.24"),
    @(44, 5, "  +4.87%  "),
    @(45, 5, "  +7.28%  "),
    @(46, 4, "'10.02"),
    @(46, 5, "  +11.45%  "),
    @(47, 4, "'16.42"),
    @(47, 5, "  -0.30%  "),
    @(48, 5, "  +1.77%  "),
    @(49, 4, "2.458.09"),
    @(49, 5, "  +0.87%  "),
    @(50, 4, "'71.39"),
    @(50, 5, "  +6.58%  "),
    @(51, 4, "'1.47"),
    @(51, 5, "  +14.89%  ")
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

Write-Output "Applied $($updates.Count) cell updates"
